# Add the new "I0" and "IF" columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - same style (s="1") as the other headers.
# Copy H1's format onto I1:J1 (paste formats only, keeps the new values).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-16
$data = @(
    @(1, 2),   # row 2
    @(5, 5),   # row 3
    @(8, 8),   # row 4
    @(1, 2),   # row 5
    @(9, 9),   # row 6
    @(8, 8),   # row 7
    @(7, 7),   # row 8
    @(8, 8),   # row 9
    @(1, 1),   # row 10
    @(9, 9),   # row 11
    @(8, 9),   # row 12
    @(5, 5),   # row 13
    @(8, 8),   # row 14
    @(1, 1),   # row 15
    @(1, 1)    # row 16
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
